$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows for target cluster "Resolving-Mac" (rows 10-13 in the original 1-based sheet)
$ws.Rows.Item(10).Resize(4).Delete()

# Update remaining rows 2-9 with recomputed TPM-derived values
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.08632499999999999
$ws.Range("H2").Value = 0.258975
$ws.Range("I2").Value = 0.04465318711422561
$ws.Range("J2").Value = 0.0446531871142256
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.4427803333333333
$ws.Range("N2").Value = 1.328341
$ws.Range("O2").Value = 0.8276247286611124
$ws.Range("P2").Value = 0.8276247286611124
$ws.Range("Q2").Value = 0.038223012275
$ws.Range("R2").Value = 0.3440071104749999
$ws.Range("S2").Value = 0.03695608186926484
$ws.Range("T2").Value = 0.03695608186926484
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.08632499999999999
$ws.Range("H3").Value = 0.258975
$ws.Range("I3").Value = 0.04465318711422561
$ws.Range("J3").Value = 0.0446531871142256
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.092221
$ws.Range("N3").Value = 0.276663
$ws.Range("O3").Value = 0.1723752713388876
$ws.Range("P3").Value = 0.1723752713388876
$ws.Range("Q3").Value = 0.007960977824999998
$ws.Range("R3").Value = 0.07164880042499999
$ws.Range("S3").Value = 0.00769710524496076
$ws.Range("T3").Value = 0.007697105244960758
$ws.Range("A4").Value = "FAPs"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.414545
$ws.Range("H4").Value = 1.243635
$ws.Range("I4").Value = 0.2144309927861761
$ws.Range("J4").Value = 0.2144309927861761
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.4427803333333333
$ws.Range("N4").Value = 1.328341
$ws.Range("O4").Value = 0.8276247286611124
$ws.Range("P4").Value = 0.8276247286611124
$ws.Range("Q4").Value = 0.1835523732816667
$ws.Range("R4").Value = 1.651971359535
$ws.Range("S4").Value = 0.177468392221192
$ws.Range("T4").Value = 0.177468392221192
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.414545
$ws.Range("H5").Value = 1.243635
$ws.Range("I5").Value = 0.2144309927861761
$ws.Range("J5").Value = 0.2144309927861761
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.092221
$ws.Range("N5").Value = 0.276663
$ws.Range("O5").Value = 0.1723752713388876
$ws.Range("P5").Value = 0.1723752713388876
$ws.Range("Q5").Value = 0.038229754445
$ws.Range("R5").Value = 0.344067790005
$ws.Range("S5").Value = 0.03696260056498417
$ws.Range("T5").Value = 0.03696260056498416
$ws.Range("A6").Value = "MuSCs"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.275124666666666
$ws.Range("H6").Value = 3.825374
$ws.Range("I6").Value = 0.6595815851101212
$ws.Range("J6").Value = 0.6595815851101212
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.4427803333333333
$ws.Range("N6").Value = 1.328341
$ws.Range("O6").Value = 0.8276247286611124
$ws.Range("P6").Value = 0.8276247286611124
$ws.Range("Q6").Value = 0.5646001249482221
$ws.Range("R6").Value = 5.081401124534
$ws.Range("S6").Value = 0.5458860304066304
$ws.Range("T6").Value = 0.5458860304066304
$ws.Range("A7").Value = "MuSCs"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.275124666666666
$ws.Range("H7").Value = 3.825374
$ws.Range("I7").Value = 0.6595815851101212
$ws.Range("J7").Value = 0.6595815851101212
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.092221
$ws.Range("N7").Value = 0.276663
$ws.Range("O7").Value = 0.1723752713388876
$ws.Range("P7").Value = 0.1723752713388876
$ws.Range("Q7").Value = 0.1175932718846666
$ws.Range("R7").Value = 1.058339446962
$ws.Range("S7").Value = 0.1136955547034908
$ws.Range("T7").Value = 0.1136955547034907
$ws.Range("A8").Value = "Resolving-Mac"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.157238
$ws.Range("H8").Value = 0.471714
$ws.Range("I8").Value = 0.08133423498947705
$ws.Range("J8").Value = 0.08133423498947705
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.4427803333333333
$ws.Range("N8").Value = 1.328341
$ws.Range("O8").Value = 0.8276247286611124
$ws.Range("P8").Value = 0.8276247286611124
$ws.Range("Q8").Value = 0.06962189405266667
$ws.Range("R8").Value = 0.6265970464739999
$ws.Range("S8").Value = 0.06731422416402509
$ws.Range("T8").Value = 0.06731422416402509
$ws.Range("A9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.157238
$ws.Range("H9").Value = 0.471714
$ws.Range("I9").Value = 0.08133423498947705
$ws.Range("J9").Value = 0.08133423498947705
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.092221
$ws.Range("N9").Value = 0.276663
$ws.Range("O9").Value = 0.1723752713388876
$ws.Range("P9").Value = 0.1723752713388876
$ws.Range("Q9").Value = 0.014500645598
$ws.Range("R9").Value = 0.130505810382
$ws.Range("S9").Value = 0.01402001082545195
$ws.Range("T9").Value = 0.01402001082545195
